$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted at row 182, pushing the
# existing rows 182-184 down to 183-185.
$ws.Rows("182:182").Insert()

# Copy the style (date format) used by column D on the row below,
# so the new date cell renders the same way.
$ws.Cells.Item(182, 4).NumberFormat = $ws.Cells.Item(183, 4).NumberFormat

$ws.Range("A182").Value = 4
$ws.Range("B182").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C182").Value = "Los Lagos"
$ws.Range("D182").Value = 44656
$ws.Range("E182").Value = 10
$ws.Range("F182").Value = "Fruta"
$ws.Range("G182").Value = 100108
$ws.Range("H182").Value = "Tropicales y subtropicales"
$ws.Range("I182").Value = 100108002
$ws.Range("J182").Value = "Mango"
$ws.Range("K182").Value = "Sin especificar"
$ws.Range("L182").Value = "Primera"
$ws.Range("M182").Value = 200
$ws.Range("N182").Value = 8000
$ws.Range("O182").Value = 8500
$ws.Range("P182").Value = 8250
$ws.Range("Q182").Value = "$/bandeja 4 kilos"
$ws.Range("R182").Value = "Perú"
$ws.Range("S182").Value = 2062
$ws.Range("T182").Value = 4
